# Regenerate the "K" column (column G) values on the active sheet.
# These values are produced upstream (Strike# -> K, recomputed std/mean
# and s_vals) and then written back into the save_data workbook, so we
# simply (re)write the recalculated K values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value, taken from the regenerated
# save_data output.
$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 2
    9  = 1
    10 = 2
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 0
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 2
    23 = 3
    24 = 1
    25 = 0
    26 = 5
    27 = 0
    28 = 1
    29 = 1
    30 = 2
    31 = 2
    32 = 2
    33 = 1
    34 = 1
    35 = 2
    36 = 0
    37 = 2
    38 = 2
    39 = 1
    40 = 2
    41 = 2
    42 = 3
    43 = 0
    44 = 1
    45 = 1
    46 = 0
    47 = 0
    48 = 0
    49 = 1
    50 = 1
    51 = 2
    52 = 2
    53 = 3
    54 = 2
    55 = 1
    56 = 1
    57 = 1
    58 = 2
    59 = 0
    60 = 0
    61 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
